# Update "想去人数" (F column) figures across the four sheets to match the
# refreshed data pull (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 2421
$ws.Cells.Item(8, 6).Value = 1788
$ws.Cells.Item(9, 6).Value = 3038
$ws.Cells.Item(11, 6).Value = 4522
$ws.Cells.Item(12, 6).Value = 403
$ws.Cells.Item(13, 6).Value = 227
$ws.Cells.Item(18, 6).Value = 251
$ws.Cells.Item(19, 6).Value = 84
$ws.Cells.Item(20, 6).Value = 117
$ws.Cells.Item(21, 6).Value = 316
$ws.Cells.Item(22, 6).Value = 4557
$ws.Cells.Item(24, 6).Value = 4099
$ws.Cells.Item(27, 6).Value = 594
$ws.Cells.Item(29, 6).Value = 96
$ws.Cells.Item(30, 6).Value = 659
$ws.Cells.Item(32, 6).Value = 603

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 35

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 211
$ws.Cells.Item(3, 6).Value = 1049

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 211
$ws.Cells.Item(4, 6).Value = 1049
$ws.Cells.Item(9, 6).Value = 2421
$ws.Cells.Item(11, 6).Value = 1788
$ws.Cells.Item(13, 6).Value = 3038
$ws.Cells.Item(15, 6).Value = 4522
$ws.Cells.Item(16, 6).Value = 403
$ws.Cells.Item(17, 6).Value = 227
$ws.Cells.Item(22, 6).Value = 251
$ws.Cells.Item(24, 6).Value = 84
$ws.Cells.Item(25, 6).Value = 117
$ws.Cells.Item(26, 6).Value = 316
$ws.Cells.Item(27, 6).Value = 4557
$ws.Cells.Item(29, 6).Value = 4099
$ws.Cells.Item(32, 6).Value = 594
$ws.Cells.Item(34, 6).Value = 97
$ws.Cells.Item(35, 6).Value = 659
$ws.Cells.Item(37, 6).Value = 603
$ws.Cells.Item(39, 6).Value = 35
